$d = $word.ActiveDocument

# 1. Insert a new paragraph "{ fat bat}" between the two blank paragraphs
#    that sit right after "Hello, {name}!" and right before "{test}".
$anchor = $d.Paragraphs.Item(7)
$anchor.Range.InsertParagraphAfter() | Out-Null
$d.Paragraphs.Item(8).Range.Text = "{ fat bat}"

# 2. Collapse the "{" / "big-test" / "}" runs (split apart by grammar-check
#    proofErr markers) in the last paragraph into a single run "{big-test}".
$d.Content.Find.Execute("{big-test}", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "{big-test}", 2) | Out-Null
